$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.030.45"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "3.148.36"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'601.91"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("D6").Value = "'142.37"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.143.60"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("D12").Value = "'0.467"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  -3.88%  "
$ws.Range("D14").Value = "'34.89"
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").Value = "3.661.13"
$ws.Range("E16").Value = "  +2.66%  "
$ws.Range("D17").Value = "64.014.91"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "3.143.74"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").Value = "'6.83"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "'487.55"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "'14.68"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "'0.712"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").Value = "'7.75"
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("D24").Value = "'88.48"
$ws.Range("E24").Value = "  +4.51%  "
$ws.Range("D25").Value = "'13.23"
$ws.Range("E25").Value = "  -5.08%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").Value = "'8.20"
$ws.Range("E28").Value = "  -7.22%  "
$ws.Range("D29").Value = "'6.97"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "'27.53"
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("D32").Value = "'0.112"
$ws.Range("E32").Value = "  -7.25%  "
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("D36").Value = "'6.06"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'52.74"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").Value = "0.0₃0743"
$ws.Range("E38").Value = "  -6.95%  "
$ws.Range("D39").Value = "'2.93"
$ws.Range("E39").Value = "  -9.07%  "
$ws.Range("D40").Value = "'0.0398"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "'431.98"
$ws.Range("E41").Value = "  -8.10%  "
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'8.38"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("D44").Value = "2.909.69"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("D45").Value = "'0.260"
$ws.Range("E45").Value = "  -4.56%  "
$ws.Range("D46").Value = "'2.19"
$ws.Range("E46").Value = "  -7.19%  "
$ws.Range("D47").Value = "'2.40"
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("D49").Value = "'25.79"
$ws.Range("E49").Value = "  -4.55%  "
$ws.Range("D50").Value = "'0.115"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "'120.85"
$ws.Range("E51").Value = "  -0.29%  "
